$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27 for the new event "610" entry
# (this shifts existing rows 27-33 down to 28-34, matching the diff)
$ws.Rows.Item(27).Insert()

# Apply the updated cell values for week 17 of 2025
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 0
$ws.Range("C4").Value = 6
$ws.Range("E4").Value = 0.16
$ws.Range("D5").Value = 12
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 0.02
$ws.Range("C8").Value = 34
$ws.Range("D8").Value = 34
$ws.Range("E8").Value = 0.07
$ws.Range("D9").Value = 2
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 0.02
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.07
$ws.Range("C13").Value = 6
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 15
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0.14
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = 0.12
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0.18
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "610"
$ws.Range("D27").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = 0.09
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0.01
$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 11
$ws.Range("E33").Value = 0.11
$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 10
$ws.Range("E34").Value = 0.12
